$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: wipe old data row, clear formatting (no fill) ---
$ws.Range("A3:O3").ClearContents()
$ws.Range("A3:O3").ClearFormats()

# --- Row 4: new "Data Type" / "Data Source" header labels (bold+underline) ---
$ws.Range("A4:O4").ClearContents()
$ws.Range("A4:O4").ClearFormats()
$ws.Range("A4").Value = "Data Type"
$ws.Range("B4").Value = "Data Source"
$ws.Range("A4:B4").Font.Bold = $true
$ws.Range("A4:B4").Font.Underline = $true

# --- Row 5: Campsite Data (keep its original fill) + reservations.csv ---
$ws.Range("C5:O5").ClearContents()
$ws.Range("C5:O5").ClearFormats()
$ws.Range("B5").ClearContents()
$ws.Range("B5").ClearFormats()
$ws.Range("A5").Value = "Campsite Data"
$ws.Range("B5").Value = "reservations.csv"

# --- Row 6: Customer Data (keep its original fill) + reservation.csv ---
$ws.Range("C6:O6").ClearContents()
$ws.Range("C6:O6").ClearFormats()
$ws.Range("B6").ClearContents()
$ws.Range("B6").ClearFormats()
$ws.Range("A6").Value = "Customer Data"
$ws.Range("B6").Value = "reservation.csv"

# --- Row 7: Census Data (keep its original fill) + census api ---
$ws.Range("A7").Value = "Census Data"
$ws.Range("B7").Value = "census api"

# --- Old rows 9,10,11 are now empty (content moved up to rows 5,6,7) ---
$ws.Range("A9").ClearContents()
$ws.Range("A9").ClearFormats()
$ws.Range("A10").ClearContents()
$ws.Range("A10").ClearFormats()
$ws.Range("A11").ClearContents()
$ws.Range("A11").ClearFormats()

# --- Column O width ---
$ws.Columns("O").ColumnWidth = 35.5

# --- Selection ---
$ws.Range("B11").Select()
